$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34/35: data for these two coins swapped position (Hedera now ranked 34, InternetComputer 35)
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0606"
$ws.Range("E34").Value = "  -0.81%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.52"
$ws.Range("E35").Value = "  -1.35%  "

# Remaining price (D) and volume (E) updates
$ws.Range("D2").Value = "'37.898.45"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "'2.036.16"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'227.48"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'60.32"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.383"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "'0.0819"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "'2.337.24"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "'14.58"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "'21.39"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("D15").Value = "'0.762"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'5.17"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "'2.012.10"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "'37.887.42"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'69.95"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'5.94"
$ws.Range("E20").Value = "  -5.11%  "
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").Value = "'224.63"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'2.43"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "'2.24"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "'166.80"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'9.31"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  -4.15%  "
$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  +5.25%  "
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D36").Value = "'6.46"
$ws.Range("E36").Value = "  +5.83%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'1.525.22"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").Value = "'17.14"
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("D42").Value = "'0.0218"
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").Value = "'96.26"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "'0.0915"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'7.11"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "'2.225.88"
$ws.Range("E51").Value = "  -0.73%  "
